$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix: slow_mvmnt msg did excess circshift in run_q.m, so events weren't
# aligned to timestamps -> the fixation/mask/prime/target duration
# probabilities (cols Q:V) were recomputed and are now uniform across all
# practice trials (rows 2-41).
for ($r = 2; $r -le 41; $r++) {
    $ws.Range("Q$r").Value = 1
    $ws.Range("R$r").Value = 0.27
    $ws.Range("S$r").Value = 0.03
    $ws.Range("T$r").Value = 0.03
    $ws.Range("U$r").Value = 0.03
    $ws.Range("V$r").Value = 0.5
}

# Update the sheet view: leftmost visible column becomes C, and the
# selection moves to S10.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S10").Select()
